# Add "AI" keyword row to keywords.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "AI"
$ws.Range("B4").Value = 1
